$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the "Urbanizacion" block (rows 14-17) formatting/structure down into the
# empty rows 18-21 so the new "Manzana" block matches the existing visual pattern.
$src = $ws.Range("B14:G17")
$dst = $ws.Range("B18:G21")
$src.Copy($dst)

# The template block (Urbanizacion) only merges columns B,C,D,E for its two
# middle rows; the new Manzana block additionally merges F and G for its
# middle rows with the same "Valor" styling used in D:E, so extend that
# formatting (and its merge) rightwards into F:G.
$srcDE = $ws.Range("D19:E20")
$dstFG = $ws.Range("F19:G20")
$srcDE.Copy($dstFG)

# Now overwrite the copied text with the new "Manzana" wording.
$ws.Range("B18").Value = "Manzana"
$ws.Range("D18").Value = "Manzana= caracteres alfanuméricos"
$ws.Range("E18").Value = "CEV<09> "
$ws.Range("F18").Value = "Manzana!= caracteres alfanuméricos"
$ws.Range("G18").Value = "CENV<13> "

$ws.Range("D19").Value = "Manzana<=1"
$ws.Range("E19").Value = "CEV<10> "
$ws.Range("F19").Value = "Manzana> 1"
$ws.Range("G19").Value = "CENV<14> "

$ws.Range("D21").Value = "Manzana= NULL"
$ws.Range("E21").Value = "CEV<11>"

# The "CEV<nn> " / "CENV<nn> " code labels elsewhere in the sheet keep their
# trailing space in a small reddish-brown run (matches existing house style),
# so reproduce that run-level formatting for the 4 new code cells.
foreach ($addr in @("E18", "G18", "E19", "G19")) {
    $cellText = $ws.Range($addr).Text
    $trailChars = $ws.Range($addr).Characters($cellText.Length, 1)
    $trailChars.Font.Size = 9
    $trailChars.Font.Color = 13209
    $trailChars.Font.Name = "Calibri"
}

# Match the new selection recorded in the workbook.
$ws.Range("B18:G21").Select()

Write-Host "done"
